$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (C) and P_Value (D) columns per corrected Diebold-Mariano test results.
$ws.Range("C2").Value = 1.09130478098523
$ws.Range("D2").Value = 0.286946905619855

$ws.Range("C3").Value = 0.9748811766525802
$ws.Range("D3").Value = 0.3402175517961528

$ws.Range("C4").Value = 0.7499397595051275
$ws.Range("D4").Value = 0.4612363096054186

$ws.Range("C5").Value = 1.178708291310626
$ws.Range("D5").Value = 0.2511061732134916

$ws.Range("C6").Value = -0.1630315300858034
$ws.Range("D6").Value = 0.87198191602162

$ws.Range("C7").Value = -0.547231061311768
$ws.Range("D7").Value = 0.5897287380262639

$ws.Range("C8").Value = 0.02062246748067501
$ws.Range("D8").Value = 0.9837327163943699

$ws.Range("C9").Value = -0.2913965006036819
$ws.Range("D9").Value = 0.7734778652850445

$ws.Range("C10").Value = 0.1562756312331871
$ws.Range("D10").Value = 0.877240947033801

$ws.Range("C11").Value = 0.4611379972359414
$ws.Range("D11").Value = 0.6492261743659635
